$wb = $excel.ActiveWorkbook

# Rename "Sheet2" -> "Negative Testing"
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Name = "Negative Testing"

# Populate the new test data (TC003 negative-testing inputs)
$ws2.Range("A1").Value = 1234454
$ws2.Range("A2").Value = '!@#$%'
$ws2.Range("A3").Value = '123sasf'
$ws2.Range("A4").Value = '!@#!123'

# Hyperlink the two "special character" entries
$ws2.Hyperlinks.Add($ws2.Range("A2"), "http://www.guru99.com", "", "", '!@#$%') | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A4"), "http://www.guru99.com", "", "", '!@#!123') | Out-Null

# A4 starts with "!" - force a quote/text prefix on it, same as the source file
$ws2.Range("A4").Value = "'!@#!123"

# Make "Negative Testing" the active sheet/tab with its own selection
$ws2.Activate() | Out-Null
$ws2.Range("A6").Select() | Out-Null
